$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.091.20"
$ws.Range("E2").Value = "  +5.01%  "
$ws.Range("D3").Value = "3.519.99"
$ws.Range("E3").Value = "  +2.90%  "
$ws.Range("E4").Value = "  +0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "594.27"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +4.13%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "169.26"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +7.07%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.521.67"
$ws.Range("E8").Value = "  +2.87%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.575"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.13%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "7.28"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.44%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.126"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +5.85%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.441"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +4.41%  "
$ws.Range("D13").Value = "4.124.59"
$ws.Range("E13").Value = "  +2.91%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "28.29"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +4.50%  "
$ws.Range("E16").Value = "  +4.52%  "
$ws.Range("D17").Value = "67.067.52"
$ws.Range("E17").Value = "  +4.87%  "
$ws.Range("D18").Value = "3.519.35"
$ws.Range("E18").Value = "  +2.18%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.33"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +4.09%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.06"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +3.31%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "396.26"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +3.60%  "
$ws.Range("E22").Value = "  +2.09%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "73.50"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +3.05%  "
$ws.Range("E24").Value = "  +11.10%  "
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("E26").Value = "  +3.35%  "
$ws.Range("E27").Value = "  +5.39%  "
$ws.Range("E28").Value = "  +2.15%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.09%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.41"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +5.52%  "
$ws.Range("E31").Value = "  +6.00%  "
$ws.Range("E32").Value = "  +4.21%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "23.65"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +3.34%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "7.48"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +7.50%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  +6.16%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "160.92"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.09%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.900"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +6.26%  "
$ws.Range("E39").Value = "  +6.10%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0753"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.23%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "4.67"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +7.08%  "
$ws.Range("E42").Value = "  +2.59%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "6.73"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +5.12%  "
$ws.Range("D44").Value = "2.844.24"
$ws.Range("E44").Value = "  +1.46%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "43.50"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.04%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "26.50"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.56"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +6.77%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0315"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.72%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "353.66"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +6.01%  "
$ws.Range("E50").Value = "  +5.03%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "33.58"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +11.82%  "
